$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly data refresh: rows 2-16 (Fecha/Volumen/Precio columns) were
# reshuffled with updated values. Apply the resulting per-cell values
# directly.
$updates = @{
    "D2"  = 44403; "J2"  = 43
    "D3"  = 44330; "J3"  = 120
    "D4"  = 44313; "J4"  = 34
    "D5"  = 44355; "J5"  = 25
    "D6"  = 44407; "J6"  = 45;  "M6"  = 5744
    "D7"  = 44341; "J7"  = 51;  "K7"  = 5500; "M7"  = 5755; "P7"  = 360
    "D8"  = 44328; "J8"  = 160; "K8"  = 6000; "M8"  = 6000; "P8"  = 375
    "D9"  = 44442; "J9"  = 25;  "L9"  = 7000; "M9"  = 6480; "P9"  = 405
    "D10" = 44350; "J10" = 25
    "D11" = 44308; "J11" = 70;  "K11" = 6000; "M11" = 6000; "P11" = 375
    "D12" = 44371; "J12" = 34;  "K12" = 5500; "M12" = 5750; "P12" = 359
    "D13" = 44363; "J13" = 160; "M13" = 5750; "P13" = 359
    "D14" = 44438;              "K14" = 5000; "M14" = 5500; "P14" = 344
    "D15" = 44306; "J15" = 50;  "K15" = 6000; "M15" = 6000; "P15" = 375
    "D16" = 44358; "J16" = 52;  "L16" = 6000; "M16" = 6000; "P16" = 375
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value2 = $updates[$addr]
}
